$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 values
$ws.Range("A6").Value = "Sprint 3"
$ws.Range("B6").Value = "open search url for search"
$ws.Range("C6").Value = "According to specs the user has the possibility to enter his own open search url. For the moment it is not possible"
$ws.Range("D6").Value = "Open"

# Match formatting used by the other data rows (left aligned, C column wraps text)
$ws.Range("A6:B6").HorizontalAlignment = -4131
$ws.Range("C6").HorizontalAlignment = -4131
$ws.Range("C6").WrapText = $true

# Row height for the new row
$ws.Rows.Item(6).RowHeight = 60

# Update the selected cell to match the saved view state
$ws.Range("C8").Select() | Out-Null
